# Natmi LR-pairs output (Col3a1-Ddr1) was regenerated with an extra sending/
# target cluster ("ECs"), growing the cross-join of clusters from 2x3 (FAPs,
# sCs) to 3x3 (ECs, FAPs, sCs) and refreshing every computed statistic.
# This rewrites data rows 2-10 (A:T) in place with the new values; Excel
# manages the shared-string table and sheet <dimension> automatically as the
# cells are written.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each row: Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
# then 16 numeric statistics columns (E..T).
$rows = @()

$rows += ,@("ECs", "Col3a1", "Ddr1", "ECs",  "2", "0.6666666666666666", "6.003292333333334", "18.009877", "0.002628468435338583", "0.002628468435338583", "3", "1", "0.2763116666666667", "0.828935", "0.02083107478128044", "0.02083107478128044", "1.658779710110556", "14.929017390995", "5.475382253677321E-05", "5.47538225367732E-05")
$rows += ,@("ECs", "Col3a1", "Ddr1", "FAPs", "2", "0.6666666666666666", "6.003292333333334", "18.009877", "0.002628468435338583", "0.002628468435338583", "3", "1", "2.180798333333333", "6.542395", "0.1644098988384798", "0.1644098988384798", "13.09196991504611", "117.827729235415", "0.0004321462295541536", "0.0004321462295541535")
$rows += ,@("ECs", "Col3a1", "Ddr1", "sCs",  "2", "0.6666666666666666", "6.003292333333334", "18.009877", "0.002628468435338583", "0.002628468435338583", "3", "1", "10.807288", "32.421864", "0.8147590263802398", "0.8147590263802398", "64.87930919452535", "583.913782750728", "0.002141568383247656", "0.002141568383247656")
$rows += ,@("FAPs", "Col3a1", "Ddr1", "ECs",  "3", "1", "2169.811686", "6509.435058", "0.9500256210433516", "0.9500256210433515", "3", "1", "0.2763116666666667", "0.828935", "0.02083107478128044", "0.02083107478128044", "599.5442833114701", "5395.89854980323", "0.01979005475608645", "0.01979005475608644")
$rows += ,@("FAPs", "Col3a1", "Ddr1", "FAPs", "3", "1", "2169.811686", "6509.435058", "0.9500256210433516", "0.9500256210433515", "3", "1", "2.180798333333333", "6.542395", "0.1644098988384798", "0.1644098988384798", "4731.921708475989", "42587.29537628391", "0.1561936162497014", "0.1561936162497013")
$rows += ,@("FAPs", "Col3a1", "Ddr1", "sCs",  "3", "1", "2169.811686", "6509.435058", "0.9500256210433516", "0.9500256210433515", "3", "1", "10.807288", "32.421864", "0.8147590263802398", "0.8147590263802398", "23449.77979636757", "211048.0181673081", "0.7740419500375638", "0.7740419500375637")
$rows += ,@("sCs", "Col3a1", "Ddr1", "ECs",  "3", "1", "108.1357256666667", "324.407177", "0.04734591052130986", "0.04734591052130985", "3", "1", "0.2763116666666667", "0.828935", "0.02083107478128044", "0.02083107478128044", "29.87916258516611", "268.912463266495", "0.000986266202657218", "0.0009862662026572178")
$rows += ,@("sCs", "Col3a1", "Ddr1", "FAPs", "3", "1", "108.1357256666667", "324.407177", "0.04734591052130986", "0.04734591052130985", "3", "1", "2.180798333333333", "6.542395", "0.1644098988384798", "0.1644098988384798", "235.8222103076572", "2122.399892768915", "0.007784136359224269", "0.007784136359224267")
$rows += ,@("sCs", "Col3a1", "Ddr1", "sCs",  "3", "1", "108.1357256666667", "324.407177", "0.04734591052130986", "0.04734591052130985", "3", "1", "10.807288", "32.421864", "0.8147590263802398", "0.8147590263802398", "1168.653930368659", "10517.88537331793", "0.03857550795942837", "0.03857550795942837")

$startRow = 2
for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowVals = $rows[$r]
    $targetRow = $startRow + $r
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $cell = $ws.Cells.Item($targetRow, $c + 1)
        if ($c -lt 4) {
            # Sending/Ligand/Receptor/Target cluster columns are text.
            $cell.Value = $rowVals[$c]
        } else {
            # Remaining columns are numeric (cast from string so scientific
            # notation literals like "5.47...E-05" parse correctly).
            $cell.Value = [double]$rowVals[$c]
        }
    }
}
